$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = 42602.584097222221
$ws.Range("B6").Value = "Bag"
$ws.Range("C6").Value = 4829
$ws.Range("D6").Value = 5621
$ws.Range("E6").Value = 656
$ws.Range("F6").Value = 95
$ws.Range("G6").Value = 41
$ws.Range("H6").Value = 69
$ws.Range("I6").Value = 29
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 99
$ws.Range("M6").Value = 0
